# Insert a new "type" column before column A, shifting the existing
# qacajobid..surveytargetdate columns (and their data) one column to the
# right, then select cell D5 to match the post-edit UI state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at A; existing columns B..N shift right and keep
# their formatting/widths/styles.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 1).Value = "type"

# Match the header formatting used by the other header cells (bold,
# centered, wrap) by copying the format from the neighboring header cell.
$headerSample = $ws.Cells.Item(1, 2)
$newHeader = $ws.Cells.Item(1, 1)
$headerSample.Copy()
$newHeader.PasteSpecial(-4122)

# Leave the selection where the author last clicked after making the edit.
[void]$ws.Range("D5").Select()
